$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "Direct read/write to Console." labels to the unified
# --- "Direct to Console." wording (column H, one row per control).
$directCells = @("H11","H12","H13","H14","H15","H16","H20","H21","H23","H24","H25","H27")
foreach ($cellRef in $directCells) {
    $ws.Range($cellRef).Value = "Direct to Console."
}

# --- Remove the "Console MVC" section (header row + its 3 data rows),
# --- which sat right after the "Guardian" row at the bottom of the table.
$ws.Rows("29:32").Delete()

# --- Column B no longer needs to fit the long "Console MVC Framework"
# --- description, so it shrinks to fit the remaining longest text.
$ws.Columns("B").ColumnWidth = 71.6

# --- Restore the active selection to the top of the table body.
$ws.Range("F3").Select()
